$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ducks")
$lo = $ws.ListObjects.Item("Table1")

# Expand table by adding 7 new rows
for ($i = 0; $i -lt 7; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# Copy formatting from the last original data row (119) into the new rows
$srcFmt = $ws.Range("A119:R119")
for ($r = 120; $r -le 126; $r++) {
    $dstFmt = $ws.Range("A" + $r + ":R" + $r)
    $srcFmt.Copy($dstFmt)
}

# Row 120
$ws.Range("A120").Value = "shark"
$ws.Range("B120").Value = 30
$ws.Range("C120").Value = "Physical Store"
$ws.Range("D120").Value = "Target"
$ws.Range("E120").Value = "Rockville"
$ws.Range("F120").Value = "MD"
$ws.Range("G120").Value = "USA"
$ws.Range("H120").Value = "USA"
$ws.Range("I120").Value = 45383
$ws.Range("J120").Value = 39.056086162518199
$ws.Range("K120").Value = -77.114897400000004
$ws.Range("L120").Value = 46
$ws.Range("M120").Value = "Allan"
$ws.Range("N120").Value = 1
$ws.Range("O120").Value = 28
$ws.Range("P120").Formula = "=2+10/16"
$ws.Range("Q120").Value = 2.25
$ws.Range("R120").Formula = "=2+15/16"

# Row 121
$ws.Range("A121").Value = "beach sunglass blue"
$ws.Range("B121").Value = 31
$ws.Range("C121").Value = "Physical Store"
$ws.Range("D121").Value = "Target"
$ws.Range("E121").Value = "Rockville"
$ws.Range("F121").Value = "MD"
$ws.Range("G121").Value = "USA"
$ws.Range("H121").Value = "USA"
$ws.Range("I121").Value = 45383
$ws.Range("J121").Value = 39.056086162518199
$ws.Range("K121").Value = -77.114897400000004
$ws.Range("L121").Value = 47
$ws.Range("M121").Value = "Allan"
$ws.Range("N121").Value = 1
$ws.Range("O121").Value = 27
$ws.Range("P121").Value = 2.25
$ws.Range("Q121").Formula = "=2+3/16"
$ws.Range("R121").Formula = "=2+6/16"

# Row 122
$ws.Range("A122").Value = "Anchors"
$ws.Range("B122").Value = 32
$ws.Range("C122").Value = "Physical Store"
$ws.Range("D122").Value = "The Boathouse"
$ws.Range("E122").Value = "Lake Buena Vista"
$ws.Range("F122").Value = "FL"
$ws.Range("G122").Value = "USA"
$ws.Range("H122").Value = "USA"
$ws.Range("I122").Value = 45407
$ws.Range("J122").Value = 28.371808918683801
$ws.Range("K122").Value = -81.517879136438097
$ws.Range("L122").Value = 48
$ws.Range("M122").Value = "Allan"
$ws.Range("N122").Value = 1
$ws.Range("O122").Value = 18
$ws.Range("P122").Value = 2
$ws.Range("Q122").Formula = "=1+10/16"
$ws.Range("R122").Value = 2.25

# Row 123
$ws.Range("A123").Value = "Octopi"
$ws.Range("B123").Value = 33
$ws.Range("C123").Value = "Claw Machine"
$ws.Range("D123").Value = "Universal Studios"
$ws.Range("E123").Value = "Orlando"
$ws.Range("F123").Value = "FL"
$ws.Range("G123").Value = "USA"
$ws.Range("H123").Value = "USA"
$ws.Range("I123").Value = 45403
$ws.Range("J123").Value = 28.4809629890139
$ws.Range("K123").Value = -81.467406691372901
$ws.Range("L123").Value = 49
$ws.Range("M123").Value = "Allan"
$ws.Range("N123").Value = 2
$ws.Range("O123").Value = 35
$ws.Range("P123").Formula = "=1+10/16"
$ws.Range("Q123").Formula = "=2+2/16"
$ws.Range("R123").Value = 2.75

# Row 124
$ws.Range("A124").Value = "Penguin"
$ws.Range("B124").Value = 34
$ws.Range("C124").Value = "Claw Machine"
$ws.Range("D124").Value = "Universal Studios"
$ws.Range("E124").Value = "Orlando"
$ws.Range("F124").Value = "FL"
$ws.Range("G124").Value = "USA"
$ws.Range("H124").Value = "USA"
$ws.Range("I124").Value = 45403
$ws.Range("J124").Value = 28.4809629890139
$ws.Range("K124").Value = -81.467406691372901
$ws.Range("L124").Value = 50
$ws.Range("M124").Value = "Allan"
$ws.Range("N124").Value = 1
$ws.Range("O124").Value = 15
$ws.Range("P124").Formula = "=2+5/16"
$ws.Range("Q124").Formula = "=2+2/16"
$ws.Range("R124").Formula = "=1+14/16"

# Row 125
$ws.Range("A125").Value = "Captain w glasses"
$ws.Range("B125").Value = 35
$ws.Range("C125").Value = "Physical Store"
$ws.Range("D125").Value = "The Boathouse"
$ws.Range("E125").Value = "Orlando"
$ws.Range("F125").Value = "FL"
$ws.Range("G125").Value = "USA"
$ws.Range("H125").Value = "USA"
$ws.Range("I125").Value = 45407
$ws.Range("J125").Value = 28.371808918683801
$ws.Range("K125").Value = -81.517879136438097
$ws.Range("L125").Value = 51
$ws.Range("M125").Value = "Allan"
$ws.Range("N125").Value = 1
$ws.Range("O125").Value = 14
$ws.Range("P125").Formula = "=2+2/16"
$ws.Range("Q125").Formula = "=1+10/16"
$ws.Range("R125").Formula = "=1+13/16"

# Row 126
$ws.Range("A126").Value = "Bucky's"
$ws.Range("B126").Value = 36
$ws.Range("C126").Value = "Physical Store"
$ws.Range("D126").Value = "Buc-ee's"
$ws.Range("E126").Value = "Lake Buena Vista"
$ws.Range("F126").Value = "FL"
$ws.Range("G126").Value = "USA"
$ws.Range("H126").Value = "USA"
$ws.Range("I126").Value = 45403
$ws.Range("J126").Value = 29.232650252891698
$ws.Range("K126").Value = -81.107395234459503
$ws.Range("L126").Value = 52
$ws.Range("M126").Value = "Derek & Cassi"
$ws.Range("N126").Value = 1
$ws.Range("O126").Value = 48
$ws.Range("P126").Value = 3
$ws.Range("Q126").Value = 2.75
$ws.Range("R126").Formula = "=3+1/16"

# Autofit columns to reflect new (longer) data
$ws.Columns.AutoFit() | Out-Null
